# noRepository-expected-generation.docx: update the bundled stack trace
# (Apache POI 4.1.0 -> 5.2.3 changed the reported JDK source line numbers)
# and keep the bold styling on the paragraph's single run.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    'ThreadPoolExecutor.runWorker(ThreadPoolExecutor.java:1130)',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'ThreadPoolExecutor.runWorker(ThreadPoolExecutor.java:1136)', 2) | Out-Null

$d.Content.Find.Execute(
    'ThreadPoolExecutor$Worker.run(ThreadPoolExecutor.java:630)',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'ThreadPoolExecutor$Worker.run(ThreadPoolExecutor.java:635)', 2) | Out-Null

$d.Content.Find.Execute(
    'java.lang.Thread.run(Thread.java:832)',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'java.lang.Thread.run(Thread.java:833)', 2) | Out-Null
